# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Naranja" at Macroferia Regional de Talca,
# right before the current row 480. This pushes the existing rows 480:520 down
# to 483:523 (dimension grows from A1:T520 to A1:T523) and populates the three
# freshly-opened rows with this week's records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 480, 481, 482 (shifts old 480:520 -> 483:523).
$ws.Range("480:482").Insert()

# --- Row 480: Fukumoto, $/bandeja 15 kilos granel ---
$ws.Range("A480").Value = 5
$ws.Range("B480").Value = "Macroferia Regional de Talca"
$ws.Range("C480").Value = "Maule"
$ws.Range("D480").Value = 44714
$ws.Range("E480").Value = 7
$ws.Range("F480").Value = "Fruta"
$ws.Range("G480").Value = 100102
$ws.Range("H480").Value = "Cítricos"
$ws.Range("I480").Value = 100102005
$ws.Range("J480").Value = "Naranja"
$ws.Range("K480").Value = "Fukumoto"
$ws.Range("L480").Value = "Primera"
$ws.Range("M480").Value = 250
$ws.Range("N480").Value = 8000
$ws.Range("O480").Value = 8000
$ws.Range("P480").Value = 8000
$ws.Range("Q480").Value = "$/bandeja 15 kilos granel"
$ws.Range("R480").Value = "Región de O'Higgins"
$ws.Range("S480").Value = 533
$ws.Range("T480").Value = 15

# --- Row 481: Fukumoto, $/bins (400 kilos) ---
$ws.Range("A481").Value = 5
$ws.Range("B481").Value = "Macroferia Regional de Talca"
$ws.Range("C481").Value = "Maule"
$ws.Range("D481").Value = 44714
$ws.Range("E481").Value = 7
$ws.Range("F481").Value = "Fruta"
$ws.Range("G481").Value = 100102
$ws.Range("H481").Value = "Cítricos"
$ws.Range("I481").Value = 100102005
$ws.Range("J481").Value = "Naranja"
$ws.Range("K481").Value = "Fukumoto"
$ws.Range("L481").Value = "Primera"
$ws.Range("M481").Value = 31
$ws.Range("N481").Value = 170000
$ws.Range("O481").Value = 200000
$ws.Range("P481").Value = 194194
$ws.Range("Q481").Value = "$/bins (400 kilos)"
$ws.Range("R481").Value = "Región de O'Higgins"
$ws.Range("S481").Value = 485
$ws.Range("T481").Value = 400

# --- Row 482: Lane Late, $/bandeja 15 kilos granel ---
$ws.Range("A482").Value = 5
$ws.Range("B482").Value = "Macroferia Regional de Talca"
$ws.Range("C482").Value = "Maule"
$ws.Range("D482").Value = 44714
$ws.Range("E482").Value = 7
$ws.Range("F482").Value = "Fruta"
$ws.Range("G482").Value = 100102
$ws.Range("H482").Value = "Cítricos"
$ws.Range("I482").Value = 100102005
$ws.Range("J482").Value = "Naranja"
$ws.Range("K482").Value = "Lane Late"
$ws.Range("L482").Value = "Primera"
$ws.Range("M482").Value = 500
$ws.Range("N482").Value = 9000
$ws.Range("O482").Value = 9000
$ws.Range("P482").Value = 9000
$ws.Range("Q482").Value = "$/bandeja 15 kilos granel"
$ws.Range("R482").Value = "Provincia de Melipilla"
$ws.Range("S482").Value = 600
$ws.Range("T482").Value = 15
